$wb = $excel.ActiveWorkbook

# Rename "Paineis DARQ" to "PAINEIS DARQ"
$wsPaineis = $wb.Worksheets.Item("Paineis DARQ")
$wsPaineis.Name = "PAINEIS DARQ"

# Rename "Recolhimento x Eliminacao" to "RECOLHIMENTO X ELIMINAÇÃO"
$wsRecolhimento = $wb.Worksheets.Item("Recolhimento x Eliminacao")
$wsRecolhimento.Name = "RECOLHIMENTO X ELIMINAÇÃO"

# Delete the "Desarquivamentos Pendentes" sheet
$wsDesarquivamentos = $wb.Worksheets.Item("Desarquivamentos Pendentes")
$excel.DisplayAlerts = $false
$wsDesarquivamentos.Delete()
$excel.DisplayAlerts = $true
